# "Add vim entry on Macro!" - append a new topic row ("macro" / "Example")
# to the bottom of the vim cheat-sheet sheet (the sheet named "Sheet2",
# which is the workbook's active/tab-selected sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Body text for the new "Example" cell (column C), exactly as authored.
$body = @'
# Intro to Macro
Vim macro is a record & paly procedure.
# Example
* Start recording: `qa` (q means start recording and a is the macro name)
* do sth. (when recording the status bar shows the word)
* Stop recording: `q`
* Invoke the macro: `5@a` (means perform the a macro from current line for 5 times)
'@

$ws.Range("A30").Value = "macro"
$ws.Range("B30").Value = "Example"
$ws.Range("C30").Value = $body

# Column C (style 14) wraps text, so Excel auto-sizes the row to fit the
# 7 wrapped lines of the new entry (7 * 15pt = 105pt), same as the other
# multi-line rows already on the sheet.
$ws.Rows.Item(30).RowHeight = 105

# Cosmetic font-substitution cleanup that came along with this save
# (the CJK UI font got swapped for Calibri on the two rich-text runs in
# the NERDTree "open recursively" note).
$rng = $ws.Range("C24")
$rng.Characters(184, 19).Font.Name = "Calibri"
$rng.Characters(203, 37).Font.Name = "Calibri"
